$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Append the two new daily GSC export rows (2025-12-05, 2025-12-06) below the
# existing data (which ended at row 60 / 2025-12-04). Leading apostrophes
# force the date-look-alike strings to be stored as literal text (matching
# the existing rows) instead of being auto-converted to date serials; the
# ClearFormats() calls drop the resulting "text" number format so the cells
# keep the sheet's default (General) style, same as every other data row.
$ws.Range("A61").Value = "'2025-12-05"
$ws.Range("B61").Value = 0
$ws.Range("C61").Value = 25

$ws.Range("A62").Value = "'2025-12-06"
$ws.Range("B62").Value = 0
$ws.Range("C62").Value = 25

$ws.Range("A61:A62").ClearFormats()
